$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.617.46"
$ws.Range("E2").Value = "  +2.02%  "

$ws.Range("D3").Value = "1.668.82"
$ws.Range("E3").Value = "  +1.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4813"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2634"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06169"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07102"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.51%  "

$ws.Range("D11").Value = "1.664.64"
$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6006"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.417"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9994"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9999"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").Value = "25.614.58"
$ws.Range("E18").Value = "  +2.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006808"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.48%  "

$ws.Range("E20").Value = "  +1.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.477"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.62%  "

$ws.Range("D22").Value = "1.878.75"
$ws.Range("E22").Value = "  +0.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.720"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.370"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.408"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "104.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.701"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.990"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.676"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07707"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04368"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9986"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.623"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6165"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9546"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.620"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8735"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9995"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01519"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.873"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3793"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.687"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1126"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.242"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05262"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.425"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3358"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.50%  "
